# DEV: Finish refactoring data-save-order service
#
# The "data" category's get/* and save/* services are being moved from the
# "selling" naming scheme to the new "orderizer" one, and the refactored
# rows (9-16) are highlighted in bold to call out the change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: data-manager-ms -> data-orders-manager ------------------------
$ws.Range("C9").Value = "data-orders-manager"
$ws.Range("E9").Value = "/orderizer/data/orders/manager/v1/"

# --- Row 10: data-get-manager-ms -> data-get-orders-manager ---------------
$ws.Range("C10").Value = "data-get-orders-manager"
$ws.Range("E10").Value = "/orderizer/data/get/orders/manager/v1/"

# --- Row 11: sales-get-ms -> data-get-search-orders ------------------------
$ws.Range("C11").Value = "data-get-search-orders"
$ws.Range("E11").Value = "/orderizer/data/get/search/orders/v1/"

# --- Row 12: sales-free-get-ms -> data-get-free-orders ---------------------
$ws.Range("C12").Value = "data-get-free-orders"
$ws.Range("E12").Value = "/orderizer/data/get/free/orders/v1/"

# --- Row 13: sales-opt-get-ms -> data-get-opt-orders ------------------------
$ws.Range("C13").Value = "data-get-opt-orders"
$ws.Range("E13").Value = "/orderizer/data/get/opt/orders/v1/"

# --- Row 14: data-save-manager-ms -> data-save-orders-manager --------------
$ws.Range("C14").Value = "data-save-orders-manager"
$ws.Range("E14").Value = "/orderizer/data/save/orders/manager/v1/"

# --- Row 15: sale-save-ms -> data-save-order --------------------------------
$ws.Range("C15").Value = "data-save-order"
$ws.Range("E15").Value = " /orderizer/data/save/order/v1/"

# --- Row 16: sales-save-ms -> data-save-orders ------------------------------
$ws.Range("C16").Value = "data-save-orders"
$ws.Range("E16").Value = " /orderizer/data/save/orders/v1/"

# --- Highlight the refactored rows in bold ---------------------------------
# Rows 10-13 keep their existing (fillId 7) fill, just gain bold text.
$ws.Range("A10:E13").Font.Bold = $true

# Rows 9, 14-16 keep their existing (fillId 5) fill, just gain bold text.
$ws.Range("A9:E9").Font.Bold = $true
$ws.Range("A14:E16").Font.Bold = $true

# --- Move the active selection to E15, matching the author's last edit -----
$ws.Range("E15").Select()
